$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("day 1")
$ws2 = $wb.Worksheets.Item("error report")

# ---------------------------------------------------------------------------
# Sheet "day 1": finish off row 106, rewrite row 107 (was the placeholder
# "Cart implementation / pending" row) and append the new Cart/Payment rows
# logged on 10th Mar,2017 (rows 108-114).
# ---------------------------------------------------------------------------

# Row 106 just needed its trailing Errors/Y-N + NA columns filled in.
$ws1.Range("G106").Value = "N"
$ws1.Range("H106").Value = "NA"

# Row 107: "Cart implementation / pending" becomes a real logged task.
$ws1.Range("B107").Value = "10th Mar,2017"
$ws1.Range("C107").Value = "Product Delete using Modal for Admin"
$ws1.Range("D107").Value = "NA"
$ws1.Range("E107").Value = "NA"
$ws1.Range("F107").Value = "60 minutes"
$ws1.Range("G107").Value = "Y"
$ws1.Range("H107").Value = "BindingResult_Error"

# D107/E107 pick up the "NA, centred / wrapped" look used elsewhere (e.g. D22).
$ws1.Range("D22").Copy()
$ws1.Range("D107").PasteSpecial(-4122)
$ws1.Range("E107").PasteSpecial(-4122)
# H107 carries the error-id style from the error-report sheet (no wrap).
$ws2.Range("A11").Copy()
$ws1.Range("H107").PasteSpecial(-4122)

# New rows for 10th Mar,2017.
$rows108to114 = @(
    @{ Row = 108; A = 107; C = "ProductDetails Page";              F = "180 minutes" },
    @{ Row = 109; A = 108; C = "Cart View Page";                   F = "60 minutes"  },
    @{ Row = 110; A = 109; C = "Cart DTO,DAO,DAOImpl";              F = "60 minutes"  },
    @{ Row = 111; A = 110; C = "CartItem DTO,DAO,DAOImpl";          F = "60 minutes"  },
    @{ Row = 112; A = 111; C = "Payment Page";                     F = "40 minutes"  },
    @{ Row = 113; A = 112; C = "Payment DTO";                      F = "120 minutes" },
    @{ Row = 114; A = 113; C = "CartFlow,CartModel,CartHandler";    F = "120 minutes" }
)

foreach ($r in $rows108to114) {
    $rowNum = $r.Row
    $ws1.Range("A$rowNum").Value = $r.A
    $ws1.Range("B$rowNum").Value = "10th Mar,2017"
    $ws1.Range("C$rowNum").Value = $r.C
    $ws1.Range("D$rowNum").Value = "NA"
    $ws1.Range("E$rowNum").Value = "NA"
    $ws1.Range("F$rowNum").Value = $r.F
    $ws1.Range("G$rowNum").Value = "N"
    $ws1.Range("H$rowNum").Value = "NA"

    $ws1.Rows.Item($rowNum).RowHeight = 28.8

    $ws1.Range("D22").Copy()
    $ws1.Range("D$rowNum").PasteSpecial(-4122)
    $ws1.Range("E$rowNum").PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# Sheet "error report": log the two new errors hit while building the
# product-delete modal / cart flow.
# ---------------------------------------------------------------------------

$ws2.Range("A15").Value = "BindingResult_Error"
$ws2.Range("B15").Value = " java.lang.IllegalStateException: Neither BindingResult nor plain target object for bean name 'product' available as request attribute"
$ws2.Range("C15").Value = "we have to return an object in the requestmapping method of controller : model.addObject(""product"",new Product()); "
$ws2.Range("D15").Value = "NA"

$chars = $ws2.Range("B15").Characters(2, 137)
$chars.Font.Underline = $true
$chars.Font.Size = 7
$chars.Font.Name = "Tahoma"
$chars.Font.Color = 0

$ws2.Range("C15").WrapText = $true
$ws2.Rows.Item(15).RowHeight = 26.4

$ws2.Range("A16").Value = "StaleStaleException_Error"
$ws2.Range("B16").Value = "Batch update returned unexpected row count from update [0]; actual row count: 0; expected: 1"
$ws2.Rows.Item(16).RowHeight = 28.8

# ---------------------------------------------------------------------------
# View state: both sheets had scrolled / selection moved on to reflect the
# newly added rows.
# ---------------------------------------------------------------------------

$ws1.Application.Goto($ws1.Range("A108"))
$ws1.Range("A108").Select()

$ws2.Application.Goto($ws2.Range("A16"))
$ws2.Range("A16").Select()
